$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row rename (column headers -> snake_case codes)
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'

# Title-case connector words (de/del/la/las/el/los/y) in state & municipality names
# plus two one-off casing fixes (GUANAJUATO -> Guanajuato, MonteMorelos -> Montemorelos)
$ws.Range("B5").Value = 'Pabellón De Arteaga'
$ws.Range("B6").Value = 'Rincón De Romos'
$ws.Range("B7").Value = 'San José De Gracia'
$ws.Range("B11").Value = 'Playas De Rosarito'
$ws.Range("B23").Value = 'Amatenango De La Frontera'
$ws.Range("B26").Value = 'Bejucal De Ocampo'
$ws.Range("B34").Value = 'Comitán De Domínguez'
$ws.Range("B49").Value = 'Ocozocoautla De Espinosa'
$ws.Range("B55").Value = 'San Cristóbal De Las Casas'
$ws.Range("B76").Value = 'Guadalupe Y Calvo'
$ws.Range("B77").Value = 'Hidalgo Del Parral'
$ws.Range("B88").Value = 'Valle De Zaragoza'
$ws.Range("B102").Value = 'San Juan De Sabinas'
$ws.Range("A112").Value = 'Ciudad De México'
$ws.Range("B116").Value = 'Cuajimalpa De Morelos'
$ws.Range("B130").Value = 'Coneto De Comonfort'
$ws.Range("B139").Value = 'Nombre De Dios'
$ws.Range("B142").Value = 'Pánuco De Coronado'
$ws.Range("B148").Value = 'San Juan De Guadalupe'
$ws.Range("B149").Value = 'San Juan Del Río'
$ws.Range("A157").Value = 'Estado De México'
$ws.Range("B157").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B160").Value = 'Almoloya De Alquisiras'
$ws.Range("B161").Value = 'Almoloya De Juárez'
$ws.Range("B162").Value = 'Almoloya Del Río'
$ws.Range("B167").Value = 'Atizapán De Zaragoza'
$ws.Range("B172").Value = 'Chapa De Mota'
$ws.Range("B174").Value = 'Coacalco De Berriozábal'
$ws.Range("B179").Value = 'Ecatepec De Morelos'
$ws.Range("B186").Value = 'Ixtapan De La Sal'
$ws.Range("B187").Value = 'Ixtapan Del Oro'
$ws.Range("B195").Value = 'Naucalpan De Juárez'
$ws.Range("B201").Value = 'San Felipe Del Progreso'
$ws.Range("B210").Value = 'Tenango Del Valle'
$ws.Range("B219").Value = 'Tlalnepantla De Baz'
$ws.Range("B224").Value = 'Valle De Bravo'
$ws.Range("B225").Value = 'Valle De Chalco Solidaridad'
$ws.Range("B226").Value = 'Villa De Allende'
$ws.Range("B227").Value = 'Villa Del Carbón'
$ws.Range("B236").Value = 'Apaseo El Alto'
$ws.Range("B237").Value = 'Apaseo El Grande'
$ws.Range("B244").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B248").Value = 'Jaral Del Progreso'
$ws.Range("B255").Value = 'Purísima Del Rincón'
$ws.Range("B260").Value = 'San Francisco Del Rincón'
$ws.Range("B262").Value = 'San Luis De La Paz'
$ws.Range("B264").Value = 'Santa Cruz De Juventino Rosas'
$ws.Range("B266").Value = 'Silao De La Victoria'
$ws.Range("B270").Value = 'Valle De Santiago'
$ws.Range("B275").Value = 'Acapulco De Juárez'
$ws.Range("B276").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B277").Value = 'Alcozauca De Guerero'
$ws.Range("B280").Value = 'Atlamajalcingo Del Monte'
$ws.Range("B281").Value = 'Atoyac De Álvarez'
$ws.Range("B282").Value = 'Ayutla De Los Libres'
$ws.Range("B284").Value = 'Buenavista De Cuéllar'
$ws.Range("B285").Value = 'Chilapa De Álvarez'
$ws.Range("B286").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B287").Value = 'Coahuayutla De José María Izazaga'
$ws.Range("B290").Value = 'Coyuca De Benítez'
$ws.Range("B291").Value = 'Coyuca De Catalán'
$ws.Range("B293").Value = 'Cuetzala Del Progreso'
$ws.Range("B294").Value = 'Cutzamala De Pinzón'
$ws.Range("B299").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B300").Value = 'Iguala De La Independencia'
$ws.Range("B301").Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range("B302").Value = 'Zihuatanejo De Azueta'
$ws.Range("B304").Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range("B315").Value = 'Taxco De Alarcón'
$ws.Range("B317").Value = 'Técpan De Galeana'
$ws.Range("B319").Value = 'Tepecoacuilco De Trujano'
$ws.Range("B320").Value = 'Tixtla De Guerero'
$ws.Range("B324").Value = 'Tlapa De Comonfort'
$ws.Range("B335").Value = 'Atotonilco El Grande'
$ws.Range("B338").Value = 'Cuautepec De Hinojosa'
$ws.Range("B343").Value = 'Huejutla De Reyes'
$ws.Range("B348").Value = 'Mixquiahuala De Juárez'
$ws.Range("B350").Value = 'Pachuca De Soto'
$ws.Range("B351").Value = 'Progreso De Obregón'
$ws.Range("B353").Value = 'Santiago De Anaya'
$ws.Range("B357").Value = 'Tepehuacán De Guerero'
$ws.Range("B358").Value = 'Tepeji Del Río De Ocampo'
$ws.Range("B359").Value = 'Tezontepec De Aldama'
$ws.Range("B363").Value = 'Tula De Allende'
$ws.Range("B364").Value = 'Tulancingo De Bravo'
$ws.Range("B367").Value = 'Zacualtipán De Ángeles'
$ws.Range("B368").Value = 'Zapotlán De Juárez'
$ws.Range("B371").Value = 'Ahualulco De Mercado'
$ws.Range("B374").Value = 'Atemajac De Brizuela'
$ws.Range("B375").Value = 'Atotonilco El Alto'
$ws.Range("B376").Value = 'Autlán De Navarro'
$ws.Range("B389").Value = 'Encarnación De Díaz'
$ws.Range("B393").Value = 'Ixtlahuacán De Los Membrillos'
$ws.Range("B394").Value = 'Ixtlahuacán Del Río'
$ws.Range("B398").Value = 'Jilotlán De Los Dolores'
$ws.Range("B402").Value = 'La Manzanilla De La Paz'
$ws.Range("B403").Value = 'Lagos De Moreno'
$ws.Range("B409").Value = 'Ojuelos De Jalisco'
$ws.Range("B413").Value = 'San Diego De Alejandría'
$ws.Range("B415").Value = 'San Juan De Los Lagos'
$ws.Range("B416").Value = 'San Juanito De Escobedo'
$ws.Range("B417").Value = 'Santa María De Los Ángeles'
$ws.Range("B418").Value = 'Santa María Del Oro'
$ws.Range("B421").Value = 'Talpa De Allende'
$ws.Range("B422").Value = 'Tamazula De Gordiano'
$ws.Range("B425").Value = 'Teocuitatlán De Corona'
$ws.Range("B426").Value = 'Tepatitlán De Morelos'
$ws.Range("B429").Value = 'Tizapán El Alto'
$ws.Range("B436").Value = 'Unión De San Antonio'
$ws.Range("B437").Value = 'Valle De Guadalupe'
$ws.Range("B438").Value = 'Valle De Juárez'
$ws.Range("B442").Value = 'Yahualica De González Gallo'
$ws.Range("B443").Value = 'Zacoalco De Torres'
$ws.Range("B446").Value = 'Zapotlán El Grande'
$ws.Range("B531").Value = 'Coatlán Del Río'
$ws.Range("B539").Value = 'Jonacatepec De Leandro Valle'
$ws.Range("B542").Value = 'Puente De Ixtla'
$ws.Range("B546").Value = 'Tetela Del Volcán'
$ws.Range("B547").Value = 'Tlaltizapán De Zapata'
$ws.Range("B558").Value = 'Ixtlán Del Río'
$ws.Range("B562").Value = 'Santa María Del Oro'
$ws.Range("B580").Value = 'San Nicolás De Los Garza'
$ws.Range("B583").Value = 'Acatlán De Pérez Figueroa'
$ws.Range("B587").Value = 'Ciénega De Zimatlán'
$ws.Range("B589").Value = 'Cuilápam De Guerero'
$ws.Range("B590").Value = 'El Barrio De La Soledad'
$ws.Range("B592").Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range("B593").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B594").Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range("B595").Value = 'Ixtlán De Juárez'
$ws.Range("B596").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B601").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B602").Value = 'Oaxaca De Juárez'
$ws.Range("B603").Value = 'Ocotlán De Morelos'
$ws.Range("B604").Value = 'Putla Villa De Guerero'
$ws.Range("B605").Value = 'Rojas De Cuauhtémoc'
$ws.Range("B609").Value = 'San Agustín De Las Juntas'
$ws.Range("B625").Value = 'San Felipe Jalapa De Díaz'
$ws.Range("B654").Value = 'San Mateo Del Mar'
$ws.Range("B657").Value = 'San Miguel Del Puerto'
$ws.Range("B665").Value = 'San Pablo Villa De Mitla'
$ws.Range("B675").Value = 'San Pedro Y San Pablo Ayutla'
$ws.Range("B676").Value = 'San Pedro Y San Pablo Teposcolula'
$ws.Range("B677").Value = 'San Pedro Y San Pablo Tequixtepec'
$ws.Range("B726").Value = 'Santo Domingo De Morelos'
$ws.Range("B733").Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range("B734").Value = 'Tanetze De Zaragoza'
$ws.Range("B735").Value = 'Tataltepec De Valdés'
$ws.Range("B736").Value = 'Teococuilco De Marcos Pérez'
$ws.Range("B737").Value = 'Teotitlán De Flores Magón'
$ws.Range("B738").Value = 'Teotitlán Del Valle'
$ws.Range("B739").Value = 'Tezoatlán De Segura Y Luna'
$ws.Range("B740").Value = 'Tlacolula De Matamoros'
$ws.Range("B744").Value = 'Villa De Chilapa De Díaz'
$ws.Range("B745").Value = 'Villa De Etla'
$ws.Range("B746").Value = 'Villa De Tamazulápam Del Progreso'
$ws.Range("B747").Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range("B748").Value = 'Villa De Zaachila'
$ws.Range("B751").Value = 'Villa Sola De Vega'
$ws.Range("B753").Value = 'Zapotitlán Del Río'
$ws.Range("B755").Value = 'Zimatlán De Álvarez'
$ws.Range("B768").Value = 'Chalchicomula De Sesma'
$ws.Range("B775").Value = 'Chila De La Sal'
$ws.Range("B786").Value = 'Huehuetlán El Chico'
$ws.Range("B787").Value = 'Huehuetlán El Grande'
$ws.Range("B789").Value = 'Ixcamilpa De Guerero'
$ws.Range("B792").Value = 'Izúcar De Matamoros'
$ws.Range("B809").Value = 'San Nicolás De Los Ranchos'
$ws.Range("B816").Value = 'Tecali De Herrera'
$ws.Range("B821").Value = 'Tepanco De López'
$ws.Range("B824").Value = 'Tepexi De Rodríguez'
$ws.Range("B825").Value = 'Tetela De Ocampo'
$ws.Range("B828").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B849").Value = 'Amealco De Bonfil'
$ws.Range("B851").Value = 'Cadereyta De Montes'
$ws.Range("B854").Value = 'Landa De Matamoros'
$ws.Range("B856").Value = 'San Juan Del Río'
$ws.Range("B867").Value = 'Ciudad Del Maíz'
$ws.Range("B877").Value = 'San Ciro De Acosta'
$ws.Range("B880").Value = 'Santa María Del Río'
$ws.Range("B886").Value = 'Tanquián De Escobedo'
$ws.Range("B887").Value = 'Villa De Guadalupe'
$ws.Range("B888").Value = 'Villa De Ramos'
$ws.Range("B889").Value = 'Villa De Reyes'
$ws.Range("B917").Value = 'Nacozari De García'
$ws.Range("B932").Value = 'Jalpa De Méndez'
$ws.Range("B955").Value = 'Soto La Marina'
$ws.Range("B968").Value = 'Papalotla De Xicohténcatl'
$ws.Range("B969").Value = 'San Pablo Del Monte'
$ws.Range("B987").Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range("B990").Value = 'Amatlán De Los Reyes'
$ws.Range("B996").Value = 'Camarón De Tejeda'
$ws.Range("B1001").Value = 'Cazones De Herrera'
$ws.Range("B1015").Value = 'Cosamaloapan De Carpio'
$ws.Range("B1023").Value = 'Hueyapan De Ocampo'
$ws.Range("B1024").Value = 'Huiloapan De Cuauhtémoc'
$ws.Range("B1025").Value = 'Ignacio De La Llave'
$ws.Range("B1027").Value = 'Ixhuatlán De Madero'
$ws.Range("B1028").Value = 'Ixhuatlán Del Café'
$ws.Range("B1036").Value = 'Juchique De Ferrer'
$ws.Range("B1040").Value = 'Lerdo De Tejada'
$ws.Range("B1042").Value = 'Martínez De La Torre'
$ws.Range("B1043").Value = 'Medellín De Bravo'
$ws.Range("B1046").Value = 'Mixtla De Altamirano'
$ws.Range("B1048").Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Range("B1060").Value = 'Paso De Ovejas'
$ws.Range("B1061").Value = 'Paso Del Macho'
$ws.Range("B1063").Value = 'Poza Rica De Hidalgo'
$ws.Range("B1072").Value = 'Sayula De Alemán'
$ws.Range("B1073").Value = 'Soledad De Doblado'
$ws.Range("B1101").Value = 'Vega De Alatorre'
$ws.Range("B1117").Value = 'El Plateado De Joaquín Amaro'
$ws.Range("B1125").Value = 'Jiménez Del Teul'
$ws.Range("B1133").Value = 'Nochistlán De Mejía'
$ws.Range("B1134").Value = 'Noria De Ángeles'
$ws.Range("B1144").Value = 'Teúl De González Ortega'
$ws.Range("B1145").Value = 'Tlaltenango De Sánchez Román'
$ws.Range("B1146").Value = 'Trinidad García De La Cadena'
$ws.Range("B1149").Value = 'Villa De Cos'
$ws.Range("A233").Value = 'Guanajuato'
$ws.Range("B577").Value = 'Montemorelos'

# Floating point last-digit recalculation refresh (1 ULP)
$ws.Range("D171").Value = 0.0009381963175794536
$ws.Range("D226").Value = 0.0009381963175794536
$ws.Range("D257").Value = 0.0009381963175794536
$ws.Range("D259").Value = 0.0009381963175794536
$ws.Range("D326").Value = 0.0009381963175794536
$ws.Range("D377").Value = 0.0009381963175794536
$ws.Range("D408").Value = 0.0009381963175794536
$ws.Range("D433").Value = 0.0009381963175794536
$ws.Range("D446").Value = 0.0009381963175794536
$ws.Range("D448").Value = 0.09346780813885304
$ws.Range("D480").Value = 0.0009381963175794536
$ws.Range("D500").Value = 0.0009381963175794536
$ws.Range("D510").Value = 0.0009381963175794536
$ws.Range("D563").Value = 0.0009381963175794536
$ws.Range("D707").Value = 0.0009381963175794536
$ws.Range("D753").Value = 0.0009381963175794536
$ws.Range("D784").Value = 0.0009381963175794536
$ws.Range("D869").Value = 0.0009381963175794536
$ws.Range("D949").Value = 0.0009381963175794536
$ws.Range("D1039").Value = 0.0009381963175794536
$ws.Range("D1085").Value = 0.0009381963175794536
$ws.Range("D1101").Value = 0.0009381963175794536
$ws.Range("D1133").Value = 0.0009381963175794536
$ws.Range("D1144").Value = 0.0009381963175794536

# Remove trailing footer/metadata rows (1158:1162) and now-empty gap row 1157
$ws.Range("A1157:D1162").EntireRow.Delete()
